$wb = $excel.ActiveWorkbook

# --- Rename header cells -----------------------------------------------
# "node" sheet: A1 "id" -> "name"
$wsNode = $wb.Worksheets.Item("node")
$wsNode.Range("A1").Value = "name"

# "edge" sheet: A1 "from" -> "Orig", B1 "to" -> "Dest"
$wsEdge = $wb.Worksheets.Item("edge")
$wsEdge.Range("A1").Value = "Orig"
$wsEdge.Range("B1").Value = "Dest"

# --- View / selection state ---------------------------------------------
# "node" keeps its selection at A2 (unchanged), but it is no longer the
# active sheet - "edge" becomes active with H14 selected instead.
$wsNode.Range("A2").Select()
$wsEdge.Activate()
$wsEdge.Range("H14").Select()
